$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Abangan Sur Elementary School -> Tabing Ilog Elementary School
$ws.Range("D2").Value = "Tabing Ilog Elementary School"
$ws.Range("E2").Value = 120.948821314155
$ws.Range("F2").Value = 14.7652274561484

# Row 6: Ibayo Elementary School -> Lias Elementary School
$ws.Range("D6").Value = "Lias Elementary School"
$ws.Range("E6").Value = 120.965390844846
$ws.Range("F6").Value = 14.7627779447143

# Row 8: Old Municipal Bldg. -> Barangay Hall Nagbalon
$ws.Range("D8").Value = "Barangay Hall Nagbalon"
$ws.Range("E8").Value = 120.950788291388
$ws.Range("F8").Value = 14.7523618894178

# Row 9: Marilao Central School -> Ramcar Covered Court
$ws.Range("D9").Value = "Ramcar Covered Court"
$ws.Range("E9").Value = 120.954403339867
$ws.Range("F9").Value = 14.7646177280722

# Row 10: Barangay Hall Nagbalon -> Old Municipal Bldg.
$ws.Range("D10").Value = "Old Municipal Bldg."
$ws.Range("E10").Value = 120.948177254006
$ws.Range("F10").Value = 14.7573006861396

# Row 17: Barangay Hall Tabing Ilog -> Tabing Ilog Elementary School
$ws.Range("D17").Value = "Tabing Ilog Elementary School"
$ws.Range("E17").Value = 120.948821314155
$ws.Range("F17").Value = 14.7652274561484
